$wb = $excel.ActiveWorkbook

# Status text changed from "Ready for handoff" to "Handed back: in sync with en-US"
$wsOv = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$wsOv.Range("B2").Value = "Handed back: in sync with en-US"
$wsOv.Range("C2").Value = "Handed back: in sync with en-US"
$wsOv.Range("B3").Value = "Handed back: in sync with en-US"
$wsOv.Range("C3").Value = "Handed back: in sync with en-US"

$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"
$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"

# Populate "Latest Target File" (F) and "Latest Handback File" (G) columns
$wsZh.Range("F2").Value = "7c52a147-6b73-4930-9b08-b0b049f4e3e2.md"
$wsZh.Range("G2").Value = "7c52a147-6b73-4930-9b08-b0b049f4e3e2.698c69b311bc70b073a0eece03175770633e26df.zh-cn.xlf"
$wsZh.Range("F3").Value = "7c52a147-6b73-4930-9b08-b0b049f4e3e2.md"
$wsZh.Range("G3").Value = "7c52a147-6b73-4930-9b08-b0b049f4e3e2.698c69b311bc70b073a0eece03175770633e26df.zh-cn.xlf"

$wsDe.Range("F2").Value = "7c52a147-6b73-4930-9b08-b0b049f4e3e2.md"
$wsDe.Range("G2").Value = "7c52a147-6b73-4930-9b08-b0b049f4e3e2.698c69b311bc70b073a0eece03175770633e26df.de-de.xlf"
$wsDe.Range("F3").Value = "7c52a147-6b73-4930-9b08-b0b049f4e3e2.md"
$wsDe.Range("G3").Value = "7c52a147-6b73-4930-9b08-b0b049f4e3e2.698c69b311bc70b073a0eece03175770633e26df.de-de.xlf"

# Latest Handback DateTime (H): zh-cn and de-de got handed back at different times
$wsZh.Range("H2").Value = "2016-03-14 09:36:21"
$wsZh.Range("H3").Value = "2016-03-14 09:36:21"
$wsDe.Range("H2").Value = "2016-03-14 09:36:36"
$wsDe.Range("H3").Value = "2016-03-14 09:36:36"
